# Regenerate merged AHB files
# - Rename the "_old" header suffixes to "_FV2210"
# - Rename the "_new" header suffixes to "_FV2304"
# - Turn the data range into an Excel Table ("Table1")
# - Freeze the header row (pane split under row 1) with a selection in the frozen pane

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSuffix = "_FV2210"
$newSuffix2 = "_FV2304"

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns A-J: "<Header>_old" -> "<Header>_FV2210"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $headers[$i] + $newSuffix
}

# Column K ("diff") stays as-is

# Columns L-U: "<Header>_new" -> "<Header>_FV2304"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $headers[$i] + $newSuffix2
}

# Turn A1:U62 into an Excel Table
$tableRange = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row and set the selection in the frozen pane
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
